$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad / Changed date) for rows 2 through 89
# from 45185 (2023-09-16) to 45204 (2023-10-05)
for ($row = 2; $row -le 89; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
